$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("main")
$ws.Range("T1").EntireColumn.Delete()
$v = $ws.Range("T1").Value()
Write-Output $v
$v2 = $ws.Range("U1").Value()
Write-Output $v2
